# Update "最低票价" (minimum ticket price) values on both the "展览" and
# "全部类型" sheets, mirroring the same edit made in both places:
#   G2: 45 -> 50
#   G3: 50 -> 60

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G2").Value = 50
    $ws.Range("G3").Value = 60
}
